$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 96, pushing existing rows 96-117 down to 97-118.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new record.
$ws.Range("A96").Value = 2
$ws.Range("B96").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C96").Value = "Coquimbo"
$ws.Range("D96").Value = 44637
$ws.Range("E96").Value = 4
$ws.Range("F96").Value = 100112024
$ws.Range("G96").Value = "Choclo"
$ws.Range("H96").Value = "Choclero"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 50000
$ws.Range("K96").Value = 230
$ws.Range("L96").Value = 250
$ws.Range("M96").Value = 240
$ws.Range("N96").Value = "$/unidad"
$ws.Range("O96").Value = "Provincia de Limarí"
$ws.Range("P96").Value = 240
$ws.Range("Q96").Value = 1
$ws.Range("R96").Value = "Hortaliza"
